# Apply updated "dSF" (column F) values for the thielbar_caleb 2022 save_data sheet.
# These correspond to a repull/recalculation of the data where column F (dSF)
# diverges from column E (dS0) for a subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    6  = -2
    11 = -1
    15 = 0
    16 = 0
    17 = 2
    20 = -1
    21 = 1
    27 = 3
    33 = -3
    36 = -1
    42 = -6
    45 = 1
    48 = 4
    52 = 2
    56 = -4
    59 = 0
    60 = 3
    65 = -3
    66 = -1
    70 = 5
    71 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
